$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "60.799.77"
$ws.Range("E2").Value = "  +2.65%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "2.708.50"
$ws.Range("E3").Value = "  +1.15%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.34%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'525.33"

# Row 6: 'Solana' -> 'Solana'
$ws.Range("D6").Value = "'144.97"
$ws.Range("E6").Value = "  -0.47%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  +0.30%  "

# Row 8: 'XRP' -> 'XRP'
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  +1.50%  "

# Row 9: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range("D9").Value = "2.706.09"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10: 'Toncoin' -> 'Toncoin'
$ws.Range("D10").Value = "'6.46"
$ws.Range("E10").Value = "  +2.85%  "

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E11").Value = "  -0.64%  "

# Row 12: 'Cardano' -> 'Cardano'
$ws.Range("E12").Value = "  -0.10%  "

# Row 13: 'TRON' -> 'TRON'
$ws.Range("E13").Value = "  +2.63%  "

# Row 14: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D14").Value = "3.157.04"
$ws.Range("E14").Value = "  +0.36%  "

# Row 15: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D15").Value = "60.714.58"
$ws.Range("E15").Value = "  +2.52%  "

# Row 16: 'Avalanche' -> 'Avalanche'
$ws.Range("D16").Value = "'21.35"
$ws.Range("E16").Value = "  +0.88%  "

# Row 17: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D18").Value = "2.703.85"
$ws.Range("E18").Value = "  -0.36%  "

# Row 19: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D19").Value = "'349.68"
$ws.Range("E19").Value = "  -1.97%  "

# Row 20: 'Polkadot' -> 'Polkadot'
$ws.Range("D20").Value = "'4.52"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21: 'Chainlink' -> 'Chainlink'
$ws.Range("E21").Value = "  +0.97%  "

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range("E22").Value = "  +1.32%  "

# Row 23: 'Dai' -> 'Dai'
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.22%  "

# Row 24: 'Litecoin' -> 'Litecoin'
$ws.Range("D24").Value = "'63.70"
$ws.Range("E24").Value = "  +2.58%  "

# Row 25: 'Polygon' -> 'Polygon'
$ws.Range("D25").Value = "'0.422"
$ws.Range("E25").Value = "  -0.32%  "

# Row 26: 'Kaspa' -> 'Kaspa'
$ws.Range("E26").Value = "  +5.13%  "

# Row 27: 'Binance-PegBSC-USD' -> 'Binance-PegBSC-USD'
$ws.Range("E27").Value = "  -0.02%  "

# Row 28: 'InternetComputer(DFINITY)' -> 'PEPE'
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0820"
$ws.Range("E28").Value = "  +0.52%  "

# Row 29: 'PEPE' -> 'InternetComputer(DFINITY)'
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'7.34"
$ws.Range("E29").Value = "  +1.15%  "

# Row 30: 'Aptos' -> 'Aptos'
$ws.Range("E30").Value = "  +8.18%  "

# Row 31: 'USDe' -> 'USDe'
$ws.Range("E31").Value = "  +0.09%  "

# Row 32: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D32").Value = "'19.29"
$ws.Range("E32").Value = "  +0.73%  "

# Row 33: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("E33").Value = "  +0.59%  "

# Row 34: 'Monero' -> 'Monero'
$ws.Range("D34").Value = "'150.31"
$ws.Range("E34").Value = "  -0.31%  "

# Row 35: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D35").Value = "'4.26"
$ws.Range("E35").Value = "  +5.26%  "

# Row 36: 'ImmutableX' -> 'ImmutableX'
$ws.Range("E36").Value = "  +9.09%  "

# Row 37: 'SuiNetwork' -> 'SuiNetwork'
$ws.Range("E37").Value = "  -2.70%  "

# Row 38: 'Fetch.AI' -> 'Fetch.AI'
$ws.Range("D38").Value = "'0.884"
$ws.Range("E38").Value = "  +3.95%  "

# Row 39: 'Stacks' -> 'Stacks'
$ws.Range("E39").Value = "  +7.19%  "

# Row 40: 'OKB' -> 'OKB'
$ws.Range("D40").Value = "'36.90"
$ws.Range("E40").Value = "  +0.28%  "

# Row 41: 'Filecoin' -> 'Filecoin'
$ws.Range("E41").Value = "  -2.62%  "

# Row 42: 'Bittensor' -> 'Bittensor'
$ws.Range("D42").Value = "'283.24"
$ws.Range("E42").Value = "  -0.20%  "

# Row 43: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D43").Value = "'20.16"
$ws.Range("E43").Value = "  +1.27%  "

# Row 44: 'Stellar' -> 'FirstDigitalUSD'
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.73%  "

# Row 45: 'Mantle' -> 'Stellar'
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0992"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46: 'FirstDigitalUSD' -> 'Mantle'
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.611"
$ws.Range("E46").Value = "  -1.54%  "

# Row 47: 'Maker' -> 'Maker'
$ws.Range("D47").Value = "2.143.35"
$ws.Range("E47").Value = "  +6.20%  "

# Row 48: 'RenderToken' -> 'Hedera'
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0540"
$ws.Range("E48").Value = "  +1.01%  "

# Row 49: 'Hedera' -> 'RenderToken'
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.90"
$ws.Range("E49").Value = "  +3.37%  "

# Row 50: 'VeChain' -> 'VeChain'
$ws.Range("E50").Value = "  +1.01%  "

# Row 51: 'WhiteBITCoin' -> 'WhiteBITCoin'
$ws.Range("D51").Value = "'10.46"
$ws.Range("E51").Value = "  +1.68%  "
